$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.254.52'
$ws.Range("E2").Value = '  +1.18%  '
$ws.Range("D3").Value = '1.857.45'
$ws.Range("E3").Value = '  +1.81%  '
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.31'
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("E7").Value = '  +0.42%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3711'
$ws.Range("E8").Value = '  +0.59%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07301'
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8922'
$ws.Range("E10").Value = '  +1.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.07'
$ws.Range("E11").Value = '  +2.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07876'
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("D13").Value = '1.837.96'
$ws.Range("E13").Value = '  +0.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.414'
$ws.Range("E14").Value = '  +1.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.519'
$ws.Range("E15").Value = '  -0.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.59'
$ws.Range("E16").Value = '  +0.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008939'
$ws.Range("E18").Value = '  +1.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.74'
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("D21").Value = '27.285.70'
$ws.Range("E21").Value = '  +1.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.087'
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("E23").Value = '  +0.17%  '
$ws.Range("D24").Value = '2.077.24'
$ws.Range("E24").Value = '  +1.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.031'
$ws.Range("E25").Value = '  +9.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.78'
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.43'
$ws.Range("E27").Value = '  +0.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.049'
$ws.Range("E28").Value = '  +0.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '116.00'
$ws.Range("E29").Value = '  +0.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.049'
$ws.Range("E30").Value = '  -0.95%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08843'
$ws.Range("E31").Value = '  -0.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.138'
$ws.Range("E32").Value = '  +6.04%  '
$ws.Range("E33").Value = '  +5.61%  '
$ws.Range("E34").Value = '  +3.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.530'
$ws.Range("E35").Value = '  +2.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.696'
$ws.Range("E36").Value = '  +9.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.108'
$ws.Range("E37").Value = '  +3.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01945'
$ws.Range("E38").Value = '  +0.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05231'
$ws.Range("E39").Value = '  +0.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.954'
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5118'
$ws.Range("E42").Value = '  -0.26%  '
$ws.Range("E43").Value = '  +0.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.538'
$ws.Range("E44").Value = '  +4.86%  '
$ws.Range("E45").Value = '  -0.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.34'
$ws.Range("E46").Value = '  +1.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.002'
$ws.Range("E47").Value = '  -0.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.77'
$ws.Range("E48").Value = '  +0.92%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.647'
$ws.Range("E49").Value = '  +1.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06201'
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '65.52'
$ws.Range("E51").Value = '  +1.55%  '
